# Apply FTUX settings changes to the workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("global_settings")

# 1) Wrap the three UI-settings header strings (row 10, columns G/H/I) in brackets.
$ws.Range("G10").Value = " [unlockedSkinPowerAsInfoBox]"
$ws.Range("H10").Value = "[showContinueButtonInUnlockedSkin]"
$ws.Range("I10").Value = "[initialMapCountdownTriggeredByPlayer]"

# 2) Extend the header band on row 8 with two more styled (blank) cells, K8:L8,
#    matching the style already used across B8:J8.
$ws.Range("J8").Copy() | Out-Null
$ws.Range("K8:L8").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Row 10 grows slightly taller to fit the new bracketed text.
$ws.Rows.Item(10).RowHeight = 192.75
